# Update "想去人数" (interested-count) values in F column for rows that
# changed between two generator runs. The same underlying data is
# duplicated on both the "展览" sheet and the "全部类型" sheet, so the
# same cell updates are applied to both worksheets.

$wb = $excel.ActiveWorkbook

# Row -> new F-column value
$updates = @{
    3  = 296
    7  = 2085
    10 = 4557
    15 = 139
    19 = 3447
    21 = 551
    31 = 684
    32 = 2084
    33 = 396
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
